$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 168119

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 66935

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 28569

$ws.Range("A5").Value = 1
$ws.Range("B5").Value = 12972
